# Braga.xlsx: remove the now-unused "Concepts" sheet and switch the
# active/selected tab to "Rules" (used in the Braga presentation), where
# the "undeclared" rule's right-hand term is updated from "s;t" to "w".

$excel.DisplayAlerts = $False

$wb = $excel.ActiveWorkbook

# 1. Delete the Concepts sheet entirely.
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Delete()

# 2. Update the Rules sheet: the "undeclared" rule's rhsTerm becomes "w"
#    instead of "s;t" (and its lhsTerm label cell mirrors that text).
$rules = $wb.Worksheets.Item("Rules")
$rules.Range("F5").Value = "w"
$rules.Range("D5").Value = "undeclared = w"

# Column D widens to fit the updated label text.
$rules.Columns("D").ColumnWidth = 13.8333333333

# Move the view's selection/cursor to D6 (single cell).
$rules.Range("D6").Select()

# 3. Rules becomes the active/selected sheet (it was Compositions before).
$rules.Activate()
